$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing "Combined" stat headers to the "Health_" prefix versions
$ws.Range("B1").Value = "Health_Index"
$ws.Range("C1").Value = "Health_Zscore"
$ws.Range("D1").Value = "Health_Pvalue"

# Add the new "Combined" stat headers at the end of the table
$ws.Range("G1").Value = "Comb_Index"
$ws.Range("H1").Value = "Comb_Zscore"
$ws.Range("I1").Value = "Comb_Pvalue"

# Approximate the resulting column widths (bestFit-style autofit)
$ws.Columns.Item(2).ColumnWidth = 12.0
$ws.Columns.Item(3).ColumnWidth = 12.8
$ws.Columns.Item(4).ColumnWidth = 13.1
$ws.Columns.Item(5).ColumnWidth = 9.6
$ws.Columns.Item(6).ColumnWidth = 10.6
$ws.Columns.Item(7).ColumnWidth = 11.3
$ws.Columns.Item(8).ColumnWidth = 12.8
$ws.Columns.Item(9).ColumnWidth = 13.1

# Match the saved selection
[void]$ws.Range("E7").Select()
